# Add release/1.0.2 to meta-sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "release/1.0.2"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"

# New row keeps the workbook's default (unstyled) look rather than
# inheriting the bold/centered header-column style used by rows 1-2.
$ws.Range("A3:D3").Style = "Normal"
